$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# The sandbox's shape id/name counter is persistent per-slide (it does not
# reuse ids after a delete). Slide 9 already has shapes with ids 1,2,6,7, so
# the very next shape created would normally be given id=3 / "TextBox 2".
# The authored slide instead has id=4 / name="TextBox 3", so we create (and
# immediately discard) one throw-away shape first to advance the counter.
$bump = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$bump.Delete()

# EMU -> point conversion (1 pt = 12700 EMU), since Shapes.AddTextbox takes
# its Left/Top/Width/Height arguments in points.
$left   = 7905750 / 12700
$top    = 6531146 / 12700
$width  = 4286250 / 12700
$height = 290913 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = $false

$tr = $tb.TextFrame.TextRange
$tr.Text = "image from https://epi.yale.edu/about-epi"
$tr.Font.Name = "+mj-lt"
$tr.Font.Bold = $false
$tr.Font.Shadow = $false
$tr.ParagraphFormat.SpaceWithin = 14.25

# AutoSize can nudge the height while the font/text is being set up, so pin
# the final size back to the authored value.
$tb.Left = $left
$tb.Top = $top
$tb.Width = $width
$tb.Height = $height
